$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing todo titles with the " todo" suffix
$ws.Range("A2").Value = "First todo"
$ws.Range("A3").Value = "Second todo"
$ws.Range("A4").Value = "Third todo"

# Add the new todo rows
$ws.Range("A5").Value = "Fourth todo"
$ws.Range("A6").Value = "Fifth todo"
$ws.Range("A7").Value = "Sixth todo"

# Set the new column B width (~19.78 characters, matching the source workbook)
$ws.Columns.Item(2).ColumnWidth = 19

# Update the active selection
$ws.Range("E15").Select()
